$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("login")
$ws2 = $wb.Worksheets.Item("user_details")

$ws2.Range("A1").Value = "Firstname"
$ws2.Range("B1").Value = "LastName"
$ws2.Range("C1").Value = "PostalCode"
$ws2.Range("A2").Value = "Ree"
$ws2.Range("B2").Value = "Test"
$ws2.Range("C2").Value = 2021

$ws2.Columns.Item(1).ColumnWidth = 13
$ws2.Columns.Item(2).ColumnWidth = 11.5
$ws2.Columns.Item(3).ColumnWidth = 13.333333333333334
$ws2.Columns.Item(4).ColumnWidth = 11.666666666666666

$ws1.Range("A4").Select()
$ws2.Activate()
$ws2.Range("C2").Select()
